$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark a few more rows as "Done"
$ws.Range("B3").Value = "Done"
$ws.Range("B4").Value = "Done"
$ws.Range("B6").Value = "Done"

# Add a remark in C2, matching the vertically-centered style used elsewhere (e.g. C11)
$ws.Range("C2").Value = "you can take care of this"
$ws.Range("C2").VerticalAlignment = $ws.Range("C11").VerticalAlignment

# Move the active selection from B8 to B7
$ws.Range("B7").Select()
